{"js": "// Remove the two blank paragraphs and the \"\u00a9 2020 ...\" copyright paragraph\n// that used to sit right after the \"LOQ4073: Qu\u00edmica Geral II (Requisito\n// fraco)\" line (and right before the trailing blank / page-break\n// paragraphs that remain at the end of the document).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOQ4073...\" paragraph by its text so this keeps working even\n// if earlier content in the document shifts paragraph indices around.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(\"LOQ4073\") !== -1) {\n        anchorIndex = i;\n        break;\n    }\n}\n\nif (anchorIndex === -1) {\n    throw new Error('Could not find the \"LOQ4073\" paragraph.');\n}\n\n// The three paragraphs immediately following the anchor are:\n//   1) an empty paragraph\n//   2) an empty paragraph with a page break before it\n//   3) the \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" copyright paragraph\n// All three are deleted.\nitems[anchorIndex + 1].delete();\nitems[anchorIndex + 2].delete();\nitems[anchorIndex + 3].delete();\n\nawait context.sync();\n", "ps1": "# Remove the two blank paragraphs and the \"\u00a9 2020 ...\" copyright paragraph\n# that used to sit right after the \"LOQ4073: Qu\u00edmica Geral II (Requisito\n# fraco)\" line (and right before the trailing blank / page-break\n# paragraphs that remain at the end of the document).\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"LOQ4073\"\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find the 'LOQ4073' paragraph.\"\n}\n$startPos = $range.Start\n\n# Locate the index of the paragraph that contains the found text so this\n# keeps working even if earlier content shifts paragraph indices around.\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Start -le $startPos -and $p.Range.End -gt $startPos) {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate the anchor paragraph index.\"\n}\n\n# Delete the three paragraphs right after the anchor paragraph:\n#   1) an empty paragraph\n#   2) an empty paragraph with a page break before it\n#   3) the \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" copyright paragraph\n$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n$d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n"}
